$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.186543703079224
$ws.Range("B1").Value = 2.14165735244751
$ws.Range("C1").Value = 6.27494478225708
$ws.Range("D1").Value = 2.304931879043579
$ws.Range("E1").Value = 1.195009827613831
